$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates per the crypto price refresh diff.
# Column D sometimes holds plain numeric-looking text (e.g. "219.71"); since the
# sheet stores these as text cells, we force NumberFormat to Text before writing
# so Excel does not silently convert them to numbers, then restore the default
# "Normal" style so no stray formatting is left behind.

$ws.Range("D2").Value = "26.274.84"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").Value = "1.666.92"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5287"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.51%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.008"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("E8").Value = "  +0.37%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06366"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.93"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07839"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.523"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.63%  "
$ws.Range("D13").Value = "1.670.84"
$ws.Range("E13").Value = "  -1.77%  "
$ws.Range("D14").Value = "1.895.85"
$ws.Range("E14").Value = "  +0.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5599"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.13%  "
$ws.Range("D16").Value = "0.0₅8114"
$ws.Range("E16").Value = "  -0.97%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.65"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.04%  "
$ws.Range("D18").Value = "26.292.80"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.008"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.726"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "200.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.90%  "
$ws.Range("E22").Value = "  +0.68%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.056"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.009"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.81"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1213"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.236"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.531"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05905"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.56%  "
$ws.Range("E31").Value = "  +0.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.513"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.326"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.598"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.16%  "
$ws.Range("E35").Value = "  +0.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.821"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5801"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01613"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.964"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.18%  "
$ws.Range("D41").Value = "1.075.33"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8577"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.008"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.79"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.45%  "
$ws.Range("D45").Value = "1.805.82"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "58.52"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.17%  "
$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.013"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.13%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4411"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.84%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.054"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.25%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₈102"
$ws.Range("E50").Value = "  -5.58%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05143"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.40%  "
